$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-02-01", "17:42:23", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:42:33", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:44:18", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:44:28", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:44:39", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:44:49", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:45:00", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:45:10", "17:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 53
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    # Leading apostrophe forces text entry so the date-shaped string
    # ("2026-02-01") is stored literally instead of being parsed into a
    # date serial number, matching the existing rows above it.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
